# "Cambios consulta de CC" - add a new backlog row reporting that the
# totals aren't shown in the "Reporte de cobranzas" (collections report).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Next free row in the backlog list (existing data goes through row 36).
$newRow = 37

$ws.Cells.Item($newRow, 1).Value = "REPORTE DE COBRANZAS NO SE VISUALIZA EL TOTAL"
$ws.Cells.Item($newRow, 2).Value = "no comenzado"

# Leave the selection where the author left it after making the edit.
$ws.Range("C29").Select()
